# Neuer Eintrag 28.01.24 Gipfeli
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Einnahmen")

$ws.Range("A9").Value = "Kiosk"
$ws.Range("B9").Value = "Gipfeli"
$ws.Range("C9").Value = (Get-Date -Year 2024 -Month 1 -Day 28 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("D9").Value = 24
$ws.Range("E9").Value = "Schulz Bäckerei"

$ws.Range("E10").Select()
